$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 531 - shifts the existing rows 531:594 down to 532:595
$ws.Rows("531:531").Insert()

# Populate the newly inserted row 531 with the new record
$ws.Range("A531").Value = 6
$ws.Range("B531").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C531").Value = "Metropolitana"
$ws.Range("D531").Value = 44610
$ws.Range("E531").Value = 13
$ws.Range("F531").Value = 100112003
$ws.Range("G531").Value = "Ajo"
$ws.Range("H531").Value = "Chino"
$ws.Range("I531").Value = "Primera"
$ws.Range("J531").Value = 900
$ws.Range("K531").Value = 16000
$ws.Range("L531").Value = 17000
$ws.Range("M531").Value = 16667
$ws.Range("N531").Value = "`$/caja 10 kilos"
$ws.Range("O531").Value = "China"
$ws.Range("P531").Value = 1667
$ws.Range("Q531").Value = 10
$ws.Range("R531").Value = "Hortaliza"
